$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 136, shifting existing rows 136-138 down to 137-139.
$ws.Rows.Item(136).Insert()

# Fill the newly inserted row 136 with the latest weekly price record.
$ws.Cells.Item(136, 1).Value = 10
$ws.Cells.Item(136, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(136, 3).Value = "La Araucanía"
$ws.Cells.Item(136, 4).Value = 45239
$ws.Cells.Item(136, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(136, 5).Value = 9
$ws.Cells.Item(136, 6).Value = 300000001
$ws.Cells.Item(136, 7).Value = "Rabanito"
$ws.Cells.Item(136, 8).Value = "Sin especificar"
$ws.Cells.Item(136, 9).Value = "Primera"
$ws.Cells.Item(136, 10).Value = 30
$ws.Cells.Item(136, 11).Value = 9000
$ws.Cells.Item(136, 12).Value = 9000
$ws.Cells.Item(136, 13).Value = 9000
$ws.Cells.Item(136, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(136, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(136, 16).Value = 750
$ws.Cells.Item(136, 17).Value = 12
$ws.Cells.Item(136, 18).Value = "Hortaliza"

# Row 137 (previously row 136) keeps its values except the region correction.
$ws.Cells.Item(137, 15).Value = "Provincia de Cautín"
